$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.061583738586743
$ws.Range("D2").Value = 1.060890972402392
$ws.Range("E2").Value = 1.067167889103512
$ws.Range("F2").Value = 1.07781382575214
$ws.Range("I2").Value = 1.051058098285538
$ws.Range("J2").Value = 1.066558667851675
$ws.Range("K2").Value = 1.063616520042916
$ws.Range("L2").Value = 1.069876477279835
$ws.Range("M2").Value = 1.080494137962325
$ws.Range("N2").Value = 1.068073303496013
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.06281500893619
$ws.Range("D3").Value = 1.061855099938356
$ws.Range("E3").Value = 1.068291085950372
$ws.Range("F3").Value = 1.079077899173472
$ws.Range("I3").Value = 1.051450017423072
$ws.Range("J3").Value = 1.067442694377056
$ws.Range("K3").Value = 1.064394784157539
$ws.Range("L3").Value = 1.070814641586134
$ws.Range("M3").Value = 1.081574887146536
$ws.Range("N3").Value = 1.068958585440459
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.063611401106162
$ws.Range("D4").Value = 1.062478591432483
$ws.Range("E4").Value = 1.069017866904393
$ws.Range("F4").Value = 1.07989610137531
$ws.Range("I4").Value = 1.051702248441203
$ws.Range("J4").Value = 1.068013854960822
$ws.Range("K4").Value = 1.064897388507492
$ws.Range("L4").Value = 1.071421103734041
$ws.Range("M4").Value = 1.082273901098576
$ws.Range("N4").Value = 1.069530557137768
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.063946129449491
$ws.Range("D5").Value = 1.06274062125233
$ws.Range("E5").Value = 1.069323406403561
$ws.Range("F5").Value = 1.080240138055887
$ws.Range("I5").Value = 1.051807959756388
$ws.Range("J5").Value = 1.068253765417447
$ws.Range("K5").Value = 1.065108448735447
$ws.Range("L5").Value = 1.071675919784111
$ws.Range("M5").Value = 1.082567695526514
$ws.Range("N5").Value = 1.06977080829476
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.064002327480994
$ws.Range("D6").Value = 1.062784612198238
$ws.Range("E6").Value = 1.069374707916324
$ws.Range("F6").Value = 1.080297907158955
$ws.Range("I6").Value = 1.051825690017317
$ws.Range("J6").Value = 1.068294035429178
$ws.Range("K6").Value = 1.065143872948008
$ws.Range("L6").Value = 1.071718696310996
$ws.Range("M6").Value = 1.082617020781951
$ws.Range("N6").Value = 1.069811135494526
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.063615874054341
$ws.Range("D7").Value = 1.062482093024431
$ws.Range("E7").Value = 1.069021949533312
$ws.Range("F7").Value = 1.07990069815479
$ws.Range("I7").Value = 1.051703662243522
$ws.Range("J7").Value = 1.068017061461732
$ws.Range("K7").Value = 1.064900209625234
$ws.Range("L7").Value = 1.071424509149555
$ws.Range("M7").Value = 1.082277827071343
$ws.Range("N7").Value = 1.069533768192277
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.06199992016792
$ws.Range("D8").Value = 1.061216879693155
$ws.Range("E8").Value = 1.067547480092631
$ws.Range("F8").Value = 1.078240972028689
$ws.Range("I8").Value = 1.051190832590105
$ws.Range("J8").Value = 1.066857608381035
$ws.Range("K8").Value = 1.063879742550241
$ws.Range("L8").Value = 1.070193657846308
$ws.Range("M8").Value = 1.080859445864701
$ws.Range("N8").Value = 1.068372668555214
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.059149838340284
$ws.Range("D9").Value = 1.058984582256707
$ws.Range("E9").Value = 1.064949194798202
$ws.Range("F9").Value = 1.075318253563852
$ws.Range("I9").Value = 1.050276664360456
$ws.Range("J9").Value = 1.064807829787901
$ws.Range("K9").Value = 1.062073960348975
$ws.Range("L9").Value = 1.068020141091113
$ws.Range("M9").Value = 1.078357688105778
$ws.Range("N9").Value = 1.066319979041365
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.057247920475428
$ws.Range("D10").Value = 1.05749441444753
$ws.Range("E10").Value = 1.063216859735621
$ws.Range("F10").Value = 1.073370962908103
$ws.Range("I10").Value = 1.049660120586649
$ws.Range("J10").Value = 1.063436738407451
$ws.Range("K10").Value = 1.060864933801399
$ws.Range("L10").Value = 1.066567953272213
$ws.Range("M10").Value = 1.076688146160929
$ws.Range("N10").Value = 1.064946940553888
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.056423895471844
$ws.Range("D11").Value = 1.056848672510943
$ws.Range("E11").Value = 1.062466684380158
$ws.Range("F11").Value = 1.07252802115725
$ws.Range("I11").Value = 1.049391457211945
$ws.Range("J11").Value = 1.062841936616644
$ws.Range("K11").Value = 1.060340169174314
$ws.Range("L11").Value = 1.065938367023043
$ws.Range("M11").Value = 1.075964789649845
$ws.Range("N11").Value = 1.064351294076313
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.056117740659442
$ws.Range("D12").Value = 1.056608740221044
$ws.Range("E12").Value = 1.06218802452295
$ws.Range("F12").Value = 1.072214950039668
$ws.Range("I12").Value = 1.049291408005703
$ws.Range("J12").Value = 1.062620832175575
$ws.Range("K12").Value = 1.06014505922632
$ws.Range("L12").Value = 1.065704391804042
$ws.Range("M12").Value = 1.075696035632219
$ws.Range("N12").Value = 1.064129875641578
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.0561834153559
$ws.Range("D13").Value = 1.056660209928185
$ws.Range("E13").Value = 1.062247798551269
$ws.Range("F13").Value = 1.072282103310723
$ws.Range("I13").Value = 1.049312880483734
$ws.Range("J13").Value = 1.062668267485651
$ws.Range("K13").Value = 1.060186919548949
$ws.Range("L13").Value = 1.065754585689567
$ws.Range("M13").Value = 1.075753687354887
$ws.Range("N13").Value = 1.064177378315235
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.056398590170667
$ws.Range("D14").Value = 1.05682884117105
$ws.Range("E14").Value = 1.062443650511468
$ws.Range("F14").Value = 1.072502141900569
$ws.Range("I14").Value = 1.049383192333986
$ws.Range("J14").Value = 1.062823663503714
$ws.Range("K14").Value = 1.060324045189269
$ws.Range("L14").Value = 1.065919028993266
$ws.Range("M14").Value = 1.07594257574153
$ws.Range("N14").Value = 1.064332995013468
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.056531156421748
$ws.Range("D15").Value = 1.056932730465866
$ws.Range("E15").Value = 1.06256431981037
$ws.Range("F15").Value = 1.072637719505199
$ws.Range("I15").Value = 1.049426479880731
$ws.Range("J15").Value = 1.06291938572832
$ws.Range("K15").Value = 1.060408507733279
$ws.Range("L15").Value = 1.06602033213779
$ws.Range("M15").Value = 1.076058947145115
$ws.Range("N15").Value = 1.064428853174612
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.057302597905048
$ws.Range("D16").Value = 1.05753725977757
$ws.Range("E16").Value = 1.063266644918353
$ws.Range("F16").Value = 1.07342691120628
$ws.Range("I16").Value = 1.04967791508493
$ws.Range("J16").Value = 1.063476189901338
$ws.Range("K16").Value = 1.060899734278096
$ws.Range("L16").Value = 1.066609720267587
$ws.Range("M16").Value = 1.076736143637303
$ws.Range("N16").Value = 1.064986448073422
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.05778637172786
$ws.Range("D17").Value = 1.057916332983692
$ws.Range("E17").Value = 1.063707176747681
$ws.Range("F17").Value = 1.073922015091818
$ws.Range("I17").Value = 1.049835179094037
$ws.Range("J17").Value = 1.063825159966629
$ws.Range("K17").Value = 1.061207532502001
$ws.Range("L17").Value = 1.066979217972176
$ws.Range("M17").Value = 1.077160813528844
$ws.Range("N17").Value = 1.065335913716233
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.058068502498055
$ws.Range("D18").Value = 1.05813739284677
$ws.Range("E18").Value = 1.063964125576172
$ws.Range("F18").Value = 1.074210824890613
$ws.Range("I18").Value = 1.049926744965076
$ws.Range("J18").Value = 1.064028601307342
$ws.Range("K18").Value = 1.061386945670837
$ws.Range("L18").Value = 1.067194664657397
$ws.Range("M18").Value = 1.077408474571343
$ws.Range("N18").Value = 1.065539643966985
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.058164694028521
$ws.Range("D19").Value = 1.058212760597553
$ws.Range("E19").Value = 1.064051737536505
$ws.Range("F19").Value = 1.074309305720813
$ws.Range("I19").Value = 1.049957938843896
$ws.Range("J19").Value = 1.064097951452457
$ws.Range("K19").Value = 1.061448100559618
$ws.Range("L19").Value = 1.067268113728442
$ws.Range("M19").Value = 1.077492913570889
$ws.Range("N19").Value = 1.065609092597261
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.057734472187355
$ws.Range("D20").Value = 1.057875666901921
$ws.Range("E20").Value = 1.063659912480069
$ws.Range("F20").Value = 1.073868892657532
$ws.Range("I20").Value = 1.049818323079366
$ws.Range("J20").Value = 1.063787729861977
$ws.Range("K20").Value = 1.061174521123134
$ws.Range("L20").Value = 1.066939582144398
$ws.Range("M20").Value = 1.077115254799551
$ws.Range("N20").Value = 1.065298430456539
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.056335228650368
$ws.Range("D21").Value = 1.056779185554622
$ws.Range("E21").Value = 1.062385977318584
$ws.Range("F21").Value = 1.072437345072811
$ws.Range("I21").Value = 1.049362494311322
$ws.Range("J21").Value = 1.062777907917291
$ws.Range("K21").Value = 1.060283670335199
$ws.Range("L21").Value = 1.065870607834125
$ws.Range("M21").Value = 1.075886954690964
$ws.Range("N21").Value = 1.064287174448863
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.05545503200859
$ws.Range("D22").Value = 1.056089349661599
$ws.Range("E22").Value = 1.061584937237555
$ws.Range("F22").Value = 1.07153747554689
$ws.Range("I22").Value = 1.049074417025903
$ws.Range("J22").Value = 1.062142016587604
$ws.Range("K22").Value = 1.059722462853604
$ws.Range("L22").Value = 1.065197812728553
$ws.Range("M22").Value = 1.075114283846183
$ws.Range("N22").Value = 1.063650380080551
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.05592168323893
$ws.Range("D23").Value = 1.056455086396962
$ws.Range("E23").Value = 1.062009590674473
$ws.Range("F23").Value = 1.072014495091953
$ws.Range("I23").Value = 1.049227272794684
$ws.Range("J23").Value = 1.062479207741025
$ws.Range("K23").Value = 1.06002007377405
$ws.Range("L23").Value = 1.065554539999741
$ws.Range("M23").Value = 1.075523928907726
$ws.Range("N23").Value = 1.063988050084086
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.057757923498816
$ws.Range("D24").Value = 1.057894042299513
$ws.Range("E24").Value = 1.063681269186538
$ws.Range("F24").Value = 1.073892896324609
$ws.Range("I24").Value = 1.04982594009225
$ws.Range("J24").Value = 1.063804643245736
$ws.Range("K24").Value = 1.061189437916122
$ws.Range("L24").Value = 1.066957492100877
$ws.Range("M24").Value = 1.077135840957218
$ws.Range("N24").Value = 1.065315367859243
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.059886971026548
$ws.Range("D25").Value = 1.059562026872656
$ws.Range("E25").Value = 1.065620932754384
$ws.Range("F25").Value = 1.076073628817231
$ws.Range("I25").Value = 1.050514246887133
$ws.Range("J25").Value = 1.065338546343833
$ws.Range("K25").Value = 1.062541704783806
$ws.Range("L25").Value = 1.068582600990791
$ws.Range("M25").Value = 1.079004745549322
$ws.Range("N25").Value = 1.06685144927568
